$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new Wins/Losses/Ties columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting (bold, border, centered) used by the other header
# cells (e.g. AC1) without disturbing the values we just set.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record (Wins/Losses/Ties) for every data row (2-54).
for ($row = 2; $row -le 54; $row++) {
    $ws.Cells.Item($row, 30).Value = 71
    $ws.Cells.Item($row, 31).Value = 91
    $ws.Cells.Item($row, 32).Value = 0
}
